$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-5 with new values
$ws.Range("A2").Value = 45088.50694444445
$ws.Range("B2").Value = 19.217
$ws.Range("C2").Value = 12.901
$ws.Range("D2").Value = 4.042
$ws.Range("E2").Value = 40.812
$ws.Range("F2").Value = 32.818
$ws.Range("G2").Value = 15.123
$ws.Range("H2").Value = 47.986
$ws.Range("I2").Value = 23.269
$ws.Range("J2").Value = 9.710000000000001
$ws.Range("K2").Value = 14.67
$ws.Range("L2").Value = 16.076
$ws.Range("M2").Value = 16.742
$ws.Range("N2").Value = 4.827
$ws.Range("O2").Value = 15.038
$ws.Range("P2").Value = 20.994
$ws.Range("Q2").Value = 12.85
$ws.Range("R2").Value = 3.46
$ws.Range("S2").Value = 2.249
$ws.Range("T2").Value = 221.547
$ws.Range("U2").Value = 41.81
$ws.Range("V2").Value = 13.881
$ws.Range("W2").Value = 27.553
$ws.Range("X2").Value = 14.055
$ws.Range("Y2").Value = 3.03
$ws.Range("Z2").Value = 24.312
$ws.Range("AA2").Value = 12.261
$ws.Range("AB2").Value = 11.125
$ws.Range("AC2").Value = 13.047
$ws.Range("AD2").Value = 16.565
$ws.Range("AE2").Value = 3.456
$ws.Range("AF2").Value = 42.557
$ws.Range("AG2").Value = 7.647
$ws.Range("AH2").Value = 17.354

$ws.Range("A3").Value = 45088.51388888889
$ws.Range("B3").Value = 3.843
$ws.Range("C3").Value = 2.174
$ws.Range("D3").Value = 1.389
$ws.Range("E3").Value = 8.113
$ws.Range("F3").Value = 6.114
$ws.Range("G3").Value = 3.026
$ws.Range("H3").Value = 16.643
$ws.Range("I3").Value = 4.654
$ws.Range("J3").Value = 1.824
$ws.Range("K3").Value = 2.492
$ws.Range("L3").Value = 3.215
$ws.Range("M3").Value = 3.217
$ws.Range("N3").Value = 0.98
$ws.Range("O3").Value = 3.008
$ws.Range("P3").Value = 4.141
$ws.Range("Q3").Value = 2.903
$ws.Range("R3").Value = 1.449
$ws.Range("S3").Value = 0.694
$ws.Range("T3").Value = 38.527
$ws.Range("U3").Value = 8.768000000000001
$ws.Range("V3").Value = 2.776
$ws.Range("W3").Value = 5.531
$ws.Range("X3").Value = 2.684
$ws.Range("Y3").Value = 0.971
$ws.Range("Z3").Value = 7.64
$ws.Range("AA3").Value = 2.452
$ws.Range("AB3").Value = 2.434
$ws.Range("AC3").Value = 2.817
$ws.Range("AD3").Value = 3.182
$ws.Range("AE3").Value = 1.265
$ws.Range("AF3").Value = 15.626
$ws.Range("AG3").Value = 1.391
$ws.Range("AH3").Value = 3.474

$ws.Range("A4").Value = 45088.52083333334
$ws.Range("B4").Value = 8.167
$ws.Range("C4").Value = 5.672
$ws.Range("D4").Value = 1.045
$ws.Range("E4").Value = 17.651
$ws.Range("F4").Value = 14.142
$ws.Range("G4").Value = 6.427
$ws.Range("H4").Value = 23.78
$ws.Range("I4").Value = 9.888999999999999
$ws.Range("J4").Value = 4.202
$ws.Range("K4").Value = 6.16
$ws.Range("L4").Value = 7.082
$ws.Range("M4").Value = 7.361
$ws.Range("N4").Value = 2.05
$ws.Range("O4").Value = 6.391
$ws.Range("P4").Value = 8.923999999999999
$ws.Range("Q4").Value = 5.62
$ws.Range("R4").Value = 0.981
$ws.Range("S4").Value = 0.575
$ws.Range("T4").Value = 89.926
$ws.Range("U4").Value = 17.784
$ws.Range("V4").Value = 5.899
$ws.Range("W4").Value = 11.673
$ws.Range("X4").Value = 6.106
$ws.Range("Y4").Value = 1.237
$ws.Range("Z4").Value = 11.399
$ws.Range("AA4").Value = 5.211
$ws.Range("AB4").Value = 4.77
$ws.Range("AC4").Value = 5.583
$ws.Range("AD4").Value = 7.322
$ws.Range("AE4").Value = 0.773
$ws.Range("AF4").Value = 21.342
$ws.Range("AG4").Value = 3.216
$ws.Range("AH4").Value = 7.375

$ws.Range("A5").Value = 45088.52777777778
$ws.Range("B5").Value = 7.69
$ws.Range("C5").Value = 5.43
$ws.Range("D5").Value = 0.8100000000000001
$ws.Range("E5").Value = 16.65
$ws.Range("F5").Value = 13.4
$ws.Range("G5").Value = 6.05
$ws.Range("H5").Value = 23.6
$ws.Range("I5").Value = 9.31
$ws.Range("J5").Value = 4.02
$ws.Range("K5").Value = 5.86
$ws.Range("L5").Value = 6.69
$ws.Range("M5").Value = 6.98
$ws.Range("N5").Value = 1.93
$ws.Range("O5").Value = 6.02
$ws.Range("P5").Value = 8.460000000000001
$ws.Range("Q5").Value = 5.25
$ws.Range("R5").Value = 0.74
$ws.Range("S5").Value = 0.46
$ws.Range("T5").Value = 84.23
$ws.Range("U5").Value = 16.84
$ws.Range("V5").Value = 5.55
$ws.Range("W5").Value = 11.13
$ws.Range("X5").Value = 5.81
$ws.Range("Y5").Value = 1.09
$ws.Range("Z5").Value = 11.37
$ws.Range("AA5").Value = 4.9
$ws.Range("AB5").Value = 4.46
$ws.Range("AC5").Value = 5.22
$ws.Range("AD5").Value = 6.95
$ws.Range("AE5").Value = 0.5600000000000001
$ws.Range("AF5").Value = 21.41
$ws.Range("AG5").Value = 3.05
$ws.Range("AH5").Value = 6.94

# Remove row 6 entirely (dataset shrinks from A1:AH6 to A1:AH5)
$ws.Rows(6).Delete()

# Adjust column widths per target (raw OOXML width = ColumnWidth + 0.8333333333333334,
# so subtract that offset to land on the exact integer width the diff expects)
$ws.Columns("B:B").ColumnWidth = 7.166666666666667
$ws.Columns("C:C").ColumnWidth = 7.166666666666667
$ws.Columns("G:G").ColumnWidth = 7.166666666666667
$ws.Columns("I:I").ColumnWidth = 7.166666666666667
$ws.Columns("L:L").ColumnWidth = 7.166666666666667
$ws.Columns("M:M").ColumnWidth = 7.166666666666667
$ws.Columns("O:O").ColumnWidth = 7.166666666666667
$ws.Columns("P:P").ColumnWidth = 7.166666666666667
$ws.Columns("T:T").ColumnWidth = 8.166666666666666
$ws.Columns("V:V").ColumnWidth = 7.166666666666667
$ws.Columns("X:X").ColumnWidth = 7.166666666666667
$ws.Columns("Z:Z").ColumnWidth = 7.166666666666667
$ws.Columns("AA:AA").ColumnWidth = 7.166666666666667
$ws.Columns("AB:AB").ColumnWidth = 7.166666666666667
$ws.Columns("AC:AC").ColumnWidth = 7.166666666666667
$ws.Columns("AD:AD").ColumnWidth = 7.166666666666667
$ws.Columns("AH:AH").ColumnWidth = 7.166666666666667
